# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp (A1)
# - Chequia overtakes Austria (rows 70/71 swap country labels)
# - Trinidad y Tobago overtakes Estonia (rows 139/140 swap country labels)
# - Montserrat overtakes Islas Malvinas (rows 214/215 swap country labels)
# - Refresh the day's case/recovered/critical/death counters for the
#   affected country rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp ---------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Septiembre de 2020 a las 18:42"

# --- Country label swaps (ranking changed between the two data pulls) --
$a70 = $ws.Range("A70").Value()
$a71 = $ws.Range("A71").Value()
$ws.Range("A70").Value = $a71
$ws.Range("A71").Value = $a70

$a139 = $ws.Range("A139").Value()
$a140 = $ws.Range("A140").Value()
$ws.Range("A139").Value = $a140
$ws.Range("A140").Value = $a139

$a214 = $ws.Range("A214").Value()
$a215 = $ws.Range("A215").Value()
$ws.Range("A214").Value = $a215
$ws.Range("A215").Value = $a214

# --- Updated counters ----------------------------------------------------
$ws.Range("B4").Value = 6559509
$ws.Range("C4").Value = 10034
$ws.Range("D4").Value = 3856749
$ws.Range("E4").Value = 2507170
$ws.Range("G4").Value = 351
$ws.Range("H4").Value = 195590

$ws.Range("B5").Value = 4547402
$ws.Range("C5").Value = 84437
$ws.Range("D5").Value = 3531212
$ws.Range("E5").Value = 940035
$ws.Range("G5").Value = 1064
$ws.Range("H5").Value = 76155

$ws.Range("B6").Value = 4210556
$ws.Range("C6").Value = 11224
$ws.Range("E6").Value = 628363
$ws.Range("G6").Value = 204
$ws.Range("H6").Value = 128857

$ws.Range("B16").Value = 358138
$ws.Range("C16").Value = 2919
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 41608

$ws.Range("B29").Value = 134653
$ws.Range("C29").Value = 359
$ws.Range("D29").Value = 118638
$ws.Range("E29").Value = 6857
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 9158

$ws.Range("B32").Value = 113206
$ws.Range("C32").Value = 1040
$ws.Range("E32").Value = 11215
$ws.Range("G32").Value = 48
$ws.Range("H32").Value = 10749

$ws.Range("D54").Value = 56558
$ws.Range("E54").Value = 644

$ws.Range("B70").Value = 31759
$ws.Range("C70").Value = 723
$ws.Range("D70").Value = 20724
$ws.Range("E70").Value = 10588
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 447

$ws.Range("B71").Value = 31247
$ws.Range("C71").Value = 664
$ws.Range("D71").Value = 26043
$ws.Range("E71").Value = 4456
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 748

$ws.Range("B78").Value = 22437
$ws.Range("C78").Value = 560
$ws.Range("D78").Value = 7312
$ws.Range("E78").Value = 14906
$ws.Range("G78").Value = 7
$ws.Range("H78").Value = 219

$ws.Range("B81").Value = 20009
$ws.Range("C81").Value = 161
$ws.Range("D81").Value = 18837
$ws.Range("E81").Value = 757

$ws.Range("B86").Value = 15414
$ws.Range("C86").Value = 121
$ws.Range("D86").Value = 12896
$ws.Range("E86").Value = 1881
$ws.Range("G86").Value = 3
$ws.Range("H86").Value = 637

$ws.Range("B91").Value = 12452
$ws.Range("C91").Value = 372
$ws.Range("E91").Value = 8351
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 297

$ws.Range("B93").Value = 10860
$ws.Range("C93").Value = 156
$ws.Range("D93").Value = 6346
$ws.Range("E93").Value = 4190
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 324

$ws.Range("B107").Value = 7088
$ws.Range("C107").Value = 65
$ws.Range("D107").Value = 6397
$ws.Range("E107").Value = 567

$ws.Range("B139").Value = 2663
$ws.Range("C139").Value = 75
$ws.Range("D139").Value = 755
$ws.Range("E139").Value = 1868
$ws.Range("G139").Value = 1
$ws.Range("H139").Value = 40

$ws.Range("B140").Value = 2600
$ws.Range("C140").Value = 15
$ws.Range("D140").Value = 2223
$ws.Range("E140").Value = 313
$ws.Range("H140").Value = 64

$ws.Range("B147").Value = 2157
$ws.Range("C147").Value = 4
$ws.Range("D147").Value = 2072
$ws.Range("E147").Value = 75

$ws.Range("B161").Value = 1315
$ws.Range("C161").Value = 2
$ws.Range("E161").Value = 39

$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
